$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Regular_Section_A")
$ws.Range("D8").Value = "CS262 (Lab) [L207]"
$ws.Range("D9").Value = "CS262 (Lab) [L207]"
$ws.Range("D25").Value = "Mon 13:00-14:30 [C101], Wed 13:00-14:30 [C101]"
$ws.Range("E25").Value = "Tue 14:30-15:30 [C101]"
$ws.Range("D26").Value = "Mon 13:00-14:30 [C102], Wed 13:00-14:30 [C102]"
$ws.Range("E26").Value = "Tue 14:30-15:30 [C102]"
$ws.Range("D27").Value = "Mon 13:00-14:30 [C104], Wed 13:00-14:30 [C104]"
$ws.Range("E27").Value = "Tue 14:30-15:30 [C104]"
$ws.Range("D28").Value = "Mon 13:00-14:30 [C202], Wed 13:00-14:30 [C202]"
$ws.Range("E28").Value = "Tue 14:30-15:30 [C202]"
$ws.Range("D29").Value = "Mon 13:00-14:30 [C203], Wed 13:00-14:30 [C203]"
$ws.Range("E29").Value = "Tue 14:30-15:30 [C203]"
$ws.Range("D30").Value = "Mon 13:00-14:30 [C204], Wed 13:00-14:30 [C204]"
$ws.Range("E30").Value = "Tue 14:30-15:30 [C204]"
$ws.Range("D31").Value = "Mon 13:00-14:30 [C205], Wed 13:00-14:30 [C205]"
$ws.Range("E31").Value = "Tue 14:30-15:30 [C205]"

$ws = $wb.Worksheets.Item("Regular_Section_B")
$ws.Range("D8").Value = "CS262 (Lab) [L207]"
$ws.Range("D9").Value = "CS262 (Lab) [L207]"
$ws.Range("D25").Value = "Mon 13:00-14:30 [C101], Wed 13:00-14:30 [C101]"
$ws.Range("E25").Value = "Tue 14:30-15:30 [C101]"
$ws.Range("E26").Value = "Tue 14:30-15:30 [C102]"
$ws.Range("E27").Value = "Tue 14:30-15:30 [C104]"
$ws.Range("D28").Value = "Mon 13:00-14:30 [C202], Wed 13:00-14:30 [C202]"
$ws.Range("E28").Value = "Tue 14:30-15:30 [C202]"
$ws.Range("D29").Value = "Mon 13:00-14:30 [C203], Wed 13:00-14:30 [C203]"
$ws.Range("E29").Value = "Tue 14:30-15:30 [C203]"
$ws.Range("D30").Value = "Mon 13:00-14:30 [C204], Wed 13:00-14:30 [C204]"
$ws.Range("E30").Value = "Tue 14:30-15:30 [C204]"
$ws.Range("D31").Value = "Mon 13:00-14:30 [C205], Wed 13:00-14:30 [C205]"
$ws.Range("E31").Value = "Tue 14:30-15:30 [C205]"

$ws = $wb.Worksheets.Item("PreMid_Section_A")
$ws.Range("B8").Value = "CS262 (Lab) [L207]"
$ws.Range("B9").Value = "CS262 (Lab) [L207]"
$ws.Range("D25").Value = "Mon 13:00-14:30 [C101], Wed 13:00-14:30 [C101]"
$ws.Range("E25").Value = "Tue 14:30-15:30 [C101]"
$ws.Range("D26").Value = "Mon 13:00-14:30 [C102], Wed 13:00-14:30 [C102]"
$ws.Range("E26").Value = "Tue 14:30-15:30 [C102]"
$ws.Range("D27").Value = "Mon 13:00-14:30 [C104], Wed 13:00-14:30 [C104]"
$ws.Range("E27").Value = "Tue 14:30-15:30 [C104]"
$ws.Range("D28").Value = "Mon 13:00-14:30 [C202], Wed 13:00-14:30 [C202]"
$ws.Range("E28").Value = "Tue 14:30-15:30 [C202]"
$ws.Range("D29").Value = "Mon 13:00-14:30 [C203], Wed 13:00-14:30 [C203]"
$ws.Range("E29").Value = "Tue 14:30-15:30 [C203]"
$ws.Range("D30").Value = "Mon 13:00-14:30 [C204], Wed 13:00-14:30 [C204]"
$ws.Range("E30").Value = "Tue 14:30-15:30 [C204]"
$ws.Range("D31").Value = "Mon 13:00-14:30 [C205], Wed 13:00-14:30 [C205]"
$ws.Range("E31").Value = "Tue 14:30-15:30 [C205]"

$ws = $wb.Worksheets.Item("PreMid_Section_B")
$ws.Range("B8").Value = "CS262 (Lab) [L106]"
$ws.Range("B9").Value = "CS262 (Lab) [L106]"
$ws.Range("D25").Value = "Mon 13:00-14:30 [C101], Wed 13:00-14:30 [C101]"
$ws.Range("E25").Value = "Tue 14:30-15:30 [C101]"
$ws.Range("E26").Value = "Tue 14:30-15:30 [C102]"
$ws.Range("E27").Value = "Tue 14:30-15:30 [C104]"
$ws.Range("D28").Value = "Mon 13:00-14:30 [C202], Wed 13:00-14:30 [C202]"
$ws.Range("E28").Value = "Tue 14:30-15:30 [C202]"
$ws.Range("D29").Value = "Mon 13:00-14:30 [C203], Wed 13:00-14:30 [C203]"
$ws.Range("E29").Value = "Tue 14:30-15:30 [C203]"
$ws.Range("D30").Value = "Mon 13:00-14:30 [C204], Wed 13:00-14:30 [C204]"
$ws.Range("E30").Value = "Tue 14:30-15:30 [C204]"
$ws.Range("D31").Value = "Mon 13:00-14:30 [C205], Wed 13:00-14:30 [C205]"
$ws.Range("E31").Value = "Tue 14:30-15:30 [C205]"

$ws = $wb.Worksheets.Item("PostMid_Section_A")
$ws.Range("B8").Value = "CS262 (Lab) [L106]"
$ws.Range("B9").Value = "CS262 (Lab) [L106]"
$ws.Range("D25").Value = "Mon 13:00-14:30 [C101], Wed 13:00-14:30 [C101]"
$ws.Range("E25").Value = "Tue 14:30-15:30 [C101]"
$ws.Range("D26").Value = "Mon 13:00-14:30 [C102], Wed 13:00-14:30 [C102]"
$ws.Range("E26").Value = "Tue 14:30-15:30 [C102]"
$ws.Range("D27").Value = "Mon 13:00-14:30 [C104], Wed 13:00-14:30 [C104]"
$ws.Range("E27").Value = "Tue 14:30-15:30 [C104]"
$ws.Range("D28").Value = "Mon 13:00-14:30 [C202], Wed 13:00-14:30 [C202]"
$ws.Range("E28").Value = "Tue 14:30-15:30 [C202]"
$ws.Range("D29").Value = "Mon 13:00-14:30 [C203], Wed 13:00-14:30 [C203]"
$ws.Range("E29").Value = "Tue 14:30-15:30 [C203]"
$ws.Range("D30").Value = "Mon 13:00-14:30 [C204], Wed 13:00-14:30 [C204]"
$ws.Range("E30").Value = "Tue 14:30-15:30 [C204]"
$ws.Range("D31").Value = "Mon 13:00-14:30 [C205], Wed 13:00-14:30 [C205]"
$ws.Range("E31").Value = "Tue 14:30-15:30 [C205]"

$ws = $wb.Worksheets.Item("PostMid_Section_B")
$ws.Range("B8").Value = "CS262 (Lab) [L106]"
$ws.Range("B9").Value = "CS262 (Lab) [L106]"
$ws.Range("D25").Value = "Mon 13:00-14:30 [C101], Wed 13:00-14:30 [C101]"
$ws.Range("E25").Value = "Tue 14:30-15:30 [C101]"
$ws.Range("E26").Value = "Tue 14:30-15:30 [C102]"
$ws.Range("E27").Value = "Tue 14:30-15:30 [C104]"
$ws.Range("D28").Value = "Mon 13:00-14:30 [C202], Wed 13:00-14:30 [C202]"
$ws.Range("E28").Value = "Tue 14:30-15:30 [C202]"
$ws.Range("D29").Value = "Mon 13:00-14:30 [C203], Wed 13:00-14:30 [C203]"
$ws.Range("E29").Value = "Tue 14:30-15:30 [C203]"
$ws.Range("D30").Value = "Mon 13:00-14:30 [C204], Wed 13:00-14:30 [C204]"
$ws.Range("E30").Value = "Tue 14:30-15:30 [C204]"
$ws.Range("D31").Value = "Mon 13:00-14:30 [C205], Wed 13:00-14:30 [C205]"
$ws.Range("E31").Value = "Tue 14:30-15:30 [C205]"
